$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) A2 (merged A2:Q2) - "Project Introduction: ..." becomes rich text with a
#    bold "Project Introduction:" lead-in.
# ---------------------------------------------------------------------------
$projIntroBold = "Project Introduction:"
$projIntroRest = " K & K JEANS is a retail store specializing in branded jeans. The company places OEM orders and sells the products under its own brand name. K & K JEANS must place orders before each season begins. As a result, accurate demand forecasting is essential. By analyzing the relationship between orders and actual demand, the company aims to identify the optimal order quantity that ensures profitability and minimizes excess inventory."
$ws.Range("A2").Value = $projIntroBold + $projIntroRest

$c = $ws.Range("A2").Characters(1, $projIntroBold.Length)
$c.Font.Bold = $true
$c.Font.Size = 12
$c.Font.Name = "Aptos Narrow"

$c = $ws.Range("A2").Characters($projIntroBold.Length + 1, $projIntroRest.Length)
$c.Font.Bold = $false
$c.Font.Size = 12
$c.Font.Name = "Aptos Narrow"

# ---------------------------------------------------------------------------
# 2) A4 - "Objective: ..." becomes rich text with a bold "Objective: " lead-in.
# ---------------------------------------------------------------------------
$objBold = "Objective: "
$objRest = " determine the break-even point, assess profitability. This analysis supports effective financial budgeting and helps sustain long-term business growth."
$ws.Range("A4").Value = $objBold + $objRest

$c = $ws.Range("A4").Characters(1, $objBold.Length)
$c.Font.Bold = $true
$c.Font.Size = 12
$c.Font.Name = "Aptos Narrow"

$c = $ws.Range("A4").Characters($objBold.Length + 1, $objRest.Length)
$c.Font.Bold = $false
$c.Font.Size = 12
$c.Font.Name = "Aptos Narrow"

# ---------------------------------------------------------------------------
# 3) A30 - "Learning point: ..." becomes rich text with a bold "Learning
#    point:" lead-in, and the tail text now reads "Conditional Formating
#    Table" instead of "Formating Table".
# ---------------------------------------------------------------------------
$lpBold = "Learning point:"
$lpRest = " Vlookup, Min(), Max(), IF() functions, Data Analysis using What-If analysis with Two-way Table, Conditional Formating Table"
$ws.Range("A30").Value = $lpBold + $lpRest

$c = $ws.Range("A30").Characters(1, $lpBold.Length)
$c.Font.Bold = $true
$c.Font.Size = 11
$c.Font.Name = "Aptos Narrow"

$c = $ws.Range("A30").Characters($lpBold.Length + 1, $lpRest.Length)
$c.Font.Bold = $false
$c.Font.Size = 11
$c.Font.Name = "Aptos Narrow"

# ---------------------------------------------------------------------------
# 4) Apply currency number formatting (with border) to the what-if data table
#    result grid G17:O25, matching the rest of the bordered table.
# ---------------------------------------------------------------------------
$currencyFmt = '_([$$-409]* #,##0.00_);_([$$-409]* \(#,##0.00\);_([$$-409]* "-"??_);_(@_)'
$ws.Range("G17:O25").NumberFormat = $currencyFmt

# ---------------------------------------------------------------------------
# 5) Update the active selection / view as left by the author.
# ---------------------------------------------------------------------------
$ws.Range("A3").Select()
$ws.Range("L9").Select()

$wb.Save()
